$wb = $excel.ActiveWorkbook

# --- Sheet2: "Data" columns C/D and "User/Pass" column E ---
$ws2 = $wb.Worksheets.Item("Sheet2")

for ($i = 1; $i -le 5; $i++) {
    $ws2.Cells.Item($i, 3).Value = "Data$i"
    $ws2.Cells.Item($i, 4).Value = "Data$i"
}

$base = 11
for ($i = 1; $i -le 5; $i++) {
    $ws2.Cells.Item($base, 5).Value = "User$i"
    $ws2.Cells.Item($base + 1, 5).Value = "Pass$i"
    $base += 10
}

# --- Sheet3: User/Pass in columns A/B ---
$ws3 = $wb.Worksheets.Item("Sheet3")

for ($i = 1; $i -le 5; $i++) {
    $ws3.Cells.Item($i, 1).Value = "User$i"
    $ws3.Cells.Item($i, 2).Value = "Pass$i"
}
